$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection (also drops any scrolled topLeftCell since we just select)
$ws.Range("R14").Select() | Out-Null

# Widen column B to 44 characters (COM ColumnWidth reads ~0.8333 narrower
# than the stored OOXML width, so compensate to land exactly on 44)
$ws.Columns.Item(2).ColumnWidth = 43.166666666666664

# Add helper column C: CONCATENATE("(",A,",","'",B,"'",")") for every data row
for ($r = 2; $r -le 85; $r++) {
    $ws.Range("C$r").Formula = '=CONCATENATE("(",A' + $r + ',",","''",B' + $r + ',"''",")")'
}
